$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range('D2') '98.545.50'
Set-TextValue $ws.Range('E2') '  -0.10%  '
Set-TextValue $ws.Range('D3') '3.368.71'
Set-TextValue $ws.Range('E3') '  +0.12%  '
Set-TextValue $ws.Range('E4') '  -0.01%  '
Set-TextValue $ws.Range('D5') '258.54'
Set-TextValue $ws.Range('E5') '  -0.37%  '
Set-TextValue $ws.Range('D6') '667.43'
Set-TextValue $ws.Range('E6') '  +6.14%  '
Set-TextValue $ws.Range('D7') '1.54'
Set-TextValue $ws.Range('E7') '  +10.21%  '
Set-TextValue $ws.Range('D8') '0.464'
Set-TextValue $ws.Range('E8') '  +19.13%  '
Set-TextValue $ws.Range('D9') '1.08'
Set-TextValue $ws.Range('E9') '  +25.28%  '
Set-TextValue $ws.Range('D11') '3.366.68'
Set-TextValue $ws.Range('E11') '  +0.23%  '
Set-TextValue $ws.Range('E12') '  +6.32%  '
Set-TextValue $ws.Range('D13') '42.26'
Set-TextValue $ws.Range('E13') '  +15.85%  '
Set-TextValue $ws.Range('E14') '  +9.00%  '
Set-TextValue $ws.Range('D15') '98.938.94'
Set-TextValue $ws.Range('E15') '  +0.51%  '
Set-TextValue $ws.Range('D16') '3.996.77'
Set-TextValue $ws.Range('E16') '  +0.13%  '
Set-TextValue $ws.Range('D17') '5.65'
Set-TextValue $ws.Range('E17') '  +3.03%  '
Set-TextValue $ws.Range('D18') '3.362.30'
Set-TextValue $ws.Range('E18') '  -0.35%  '
Set-TextValue $ws.Range('D19') '7.64'
Set-TextValue $ws.Range('E19') '  +25.55%  '
Set-TextValue $ws.Range('D20') '16.85'
Set-TextValue $ws.Range('E20') '  +10.80%  '
Set-TextValue $ws.Range('B21') 'SuiNetwork'
Set-TextValue $ws.Range('C21') 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextValue $ws.Range('D21') '3.59'
Set-TextValue $ws.Range('E21') '  +1.12%  '
Set-TextValue $ws.Range('B22') 'BitcoinCash'
Set-TextValue $ws.Range('C22') 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue $ws.Range('D22') '531.44'
Set-TextValue $ws.Range('E22') '  +8.69%  '
Set-TextValue $ws.Range('D23') '10.53'
Set-TextValue $ws.Range('E23') '  +12.32%  '
Set-TextValue $ws.Range('D24') '0.0000218'
Set-TextValue $ws.Range('E24') '  +3.85%  '
Set-TextValue $ws.Range('D25') '0.440'
Set-TextValue $ws.Range('E25') '  +55.77%  '
Set-TextValue $ws.Range('D26') '102.65'
Set-TextValue $ws.Range('E26') '  +15.30%  '
Set-TextValue $ws.Range('D27') '6.25'
Set-TextValue $ws.Range('E27') '  +10.93%  '
Set-TextValue $ws.Range('D28') '12.65'
Set-TextValue $ws.Range('E28') '  +5.99%  '
Set-TextValue $ws.Range('D29') '3.545.55'
Set-TextValue $ws.Range('E29') '  -0.09%  '
Set-TextValue $ws.Range('D30') '0.150'
Set-TextValue $ws.Range('E30') '  +10.28%  '
Set-TextValue $ws.Range('E31') '  -0.42%  '
Set-TextValue $ws.Range('D32') '11.05'
Set-TextValue $ws.Range('E32') '  +14.62%  '
Set-TextValue $ws.Range('E33') '  -0.47%  '
Set-TextValue $ws.Range('D34') '0.986'
Set-TextValue $ws.Range('E34') '  -1.57%  '
Set-TextValue $ws.Range('D35') '29.61'
Set-TextValue $ws.Range('E35') '  +5.45%  '
Set-TextValue $ws.Range('D36') '0.547'
Set-TextValue $ws.Range('D37') '7.90'
Set-TextValue $ws.Range('E37') '  +8.38%  '
Set-TextValue $ws.Range('E38') '  +8.38%  '
Set-TextValue $ws.Range('D39') '0.160'
Set-TextValue $ws.Range('E39') '  +6.19%  '
Set-TextValue $ws.Range('D40') '530.65'
Set-TextValue $ws.Range('E40') '  +6.11%  '
Set-TextValue $ws.Range('E41') '  +6.98%  '
Set-TextValue $ws.Range('E42') '  -0.82%  '
Set-TextValue $ws.Range('D43') '0.0437'
Set-TextValue $ws.Range('E43') '  +33.63%  '
Set-TextValue $ws.Range('D44') '3.76'
Set-TextValue $ws.Range('E44') '  +0.90%  '
Set-TextValue $ws.Range('D45') '3.46'
Set-TextValue $ws.Range('E45') '  +5.31%  '
Set-TextValue $ws.Range('D46') '0.836'
Set-TextValue $ws.Range('E46') '  +5.10%  '
Set-TextValue $ws.Range('E47') '  +0.03%  '
Set-TextValue $ws.Range('B48') 'Stacks'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range('D48') '2.09'
Set-TextValue $ws.Range('E48') '  +7.64%  '
Set-TextValue $ws.Range('B49') 'Cosmos'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range('D49') '7.97'
Set-TextValue $ws.Range('E49') '  +19.03%  '
Set-TextValue $ws.Range('D50') '5.13'
Set-TextValue $ws.Range('E50') '  +10.89%  '
Set-TextValue $ws.Range('E51') '  +12.42%  '
